# JADE_Scrum.xlsx - Update sprint backlog (Sprint 06 Backlog progress + Product Backlog status)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Product Backlog sheet: mark the two "In Work" P&L/Reports tasks finished
# ---------------------------------------------------------------------------
$product = $wb.Worksheets.Item("Product Backlog")
$product.Range("F42").Value = "Finished in Sprint 6"
$product.Range("F43").Value = "Finished in Sprint 6"

# View state: scrolled down / selection moved to F43, no longer the active tab
$product.Range("A43").Select()
$product.Range("F43").Select()

# ---------------------------------------------------------------------------
# 2) Sprint 06 Backlog sheet: fill in task statuses + new tasks for days 5/6
# ---------------------------------------------------------------------------
$sprint6 = $wb.Worksheets.Item("Sprint 06 Backlog")

# Previously-finished tasks get their completion day recorded
$sprint6.Range("E19").Value = "Completed Day 6"
$sprint6.Range("E27").Value = "Completed Day 5"
$sprint6.Range("E28").Value = "Completed Day 5"
$sprint6.Range("E29").Value = "Completed Day 6"

# New tasks added to the backlog (rows 30-42)
$sprint6.Range("B30").Value = "P&L"
$sprint6.Range("D30").Value = "Implement PnLReport actionlistener"
$sprint6.Range("E30").Value = "Completed Day 6"

$sprint6.Range("B31").Value = "RTM"
$sprint6.Range("D31").Value = "Implement serverReport actionlistener"
$sprint6.Range("E31").Value = "Completed Day 6"

$sprint6.Range("B32").Value = "P&L"
$sprint6.Range("D32").Value = "Add income, cost, and profit methods  to Order.java"
$sprint6.Range("E32").Value = "Completed Day 5"

$sprint6.Range("B33").Value = "P&L"
$sprint6.Range("D33").Value = "Add PnLReport to Order.java"
$sprint6.Range("E33").Value = "Completed Day 5"

$sprint6.Range("B34").Value = "P&L"
$sprint6.Range("D34").Value = "Add PnLReport to Store.java"
$sprint6.Range("E34").Value = "Completed Day 5"

$sprint6.Range("B35").Value = "RTM"
$sprint6.Range("D35").Value = "Add filledOrder attribute to Server.java"
$sprint6.Range("E35").Value = "Completed Day 5"

$sprint6.Range("B36").Value = "RTM"
$sprint6.Range("D36").Value = "Add fillOrder and filledOrders methods to Server.java"
$sprint6.Range("E36").Value = "Completed Day 5"

$sprint6.Range("B37").Value = "RTM"
$sprint6.Range("D37").Value = "Add serverReport method to Server.java "
$sprint6.Range("E37").Value = "Completed Day 6"

$sprint6.Range("B38").Value = "RTM"
$sprint6.Range("D38").Value = "Add serverReport method to Store.java"
$sprint6.Range("E38").Value = "Completed Day 6"

$sprint6.Range("B39").Value = "RTM"
$sprint6.Range("D39").Value = "Update Program and file verison "
$sprint6.Range("E39").Value = "Completed Day 6"

$sprint6.Range("B40").Value = "RTM"
$sprint6.Range("D40").Value = "Add fill method to Order.java"
$sprint6.Range("E40").Value = "Completed Day 6"

$sprint6.Range("B41").Value = "RTM"
$sprint6.Range("D41").Value = "include filledOrders attributes in saving and loading from files"
$sprint6.Range("E41").Value = "Completed Day 6"

$sprint6.Range("D42").Value = "Add Report buttons to toolbar"
$sprint6.Range("E42").Value = "Completed Day 6"

# View state: this tab is now the active/selected one, scrolled to top, with E19 selected
$sprint6.Range("E19").Select()
$sprint6.Activate()
